# Refresh the cryptos list: per-row Price (D) / Volume(1h) (E) values, plus
# Chainlink / Polygon / WrappedEther re-ranked into rows 13-15.
#
# The sheet stores these columns as TEXT (original cells are inline strings,
# e.g. "230.02", "  -0.54%  "). A leading apostrophe forces Excel to keep
# numeric-looking values ("230.02", "97.90", ...) as text instead of silently
# coercing them to numbers (which would drop meaningful trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.593.59"
$ws.Range("E2").Value = "'  -2.59%  "

$ws.Range("D3").Value = "'1.806.54"
$ws.Range("E3").Value = "'  -1.84%  "

$ws.Range("E4").Value = "'  +0.51%  "

$ws.Range("D5").Value = "'230.02"
$ws.Range("E5").Value = "'  -0.54%  "

$ws.Range("E6").Value = "'  +0.00%  "

$ws.Range("E7").Value = "'  +0.47%  "

$ws.Range("D8").Value = "'39.07"
$ws.Range("E8").Value = "'  -10.90%  "

$ws.Range("E9").Value = "'  +3.98%  "

$ws.Range("D10").Value = "'0.0679"
$ws.Range("E10").Value = "'  -3.22%  "

$ws.Range("E11").Value = "'  -1.87%  "

$ws.Range("D12").Value = "'2.069.45"
$ws.Range("E12").Value = "'  -1.74%  "

$ws.Range("B13").Value = "'WrappedEther"
$ws.Range("C13").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.827.50"
$ws.Range("E13").Value = "'  -0.75%  "

$ws.Range("B14").Value = "'Chainlink"
$ws.Range("C14").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'11.16"
$ws.Range("E14").Value = "'  -0.82%  "

$ws.Range("B15").Value = "'Polygon"
$ws.Range("C15").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.661"
$ws.Range("E15").Value = "'  -2.00%  "

$ws.Range("E16").Value = "'  -3.28%  "

$ws.Range("D17").Value = "'34.607.91"
$ws.Range("E17").Value = "'  -2.46%  "

$ws.Range("D18").Value = "'69.14"
$ws.Range("E18").Value = "'  -1.66%  "

$ws.Range("D19").Value = "'242.86"
$ws.Range("E19").Value = "'  -0.63%  "

$ws.Range("D20").Value = "'0.0₃0781"
$ws.Range("E20").Value = "'  -2.47%  "

$ws.Range("D21").Value = "'11.84"
$ws.Range("E21").Value = "'  -1.85%  "

$ws.Range("D22").Value = "'4.66"
$ws.Range("E22").Value = "'  -0.99%  "

$ws.Range("E23").Value = "'  +0.57%  "

$ws.Range("D24").Value = "'2.24"
$ws.Range("E24").Value = "'  +1.14%  "

$ws.Range("D25").Value = "'171.95"
$ws.Range("E25").Value = "'  +0.08%  "

$ws.Range("E26").Value = "'  -2.48%  "

$ws.Range("D27").Value = "'17.20"
$ws.Range("E27").Value = "'  -3.38%  "

$ws.Range("E28").Value = "'  +0.07%  "

$ws.Range("D29").Value = "'1.48"
$ws.Range("E29").Value = "'  -5.95%  "

$ws.Range("E30").Value = "'  +0.50%  "

$ws.Range("E31").Value = "'  +2.82%  "

$ws.Range("D32").Value = "'0.0542"
$ws.Range("E32").Value = "'  -1.92%  "

$ws.Range("D33").Value = "'3.89"
$ws.Range("E33").Value = "'  -4.40%  "

$ws.Range("E34").Value = "'  +11.77%  "

$ws.Range("E35").Value = "'  -3.55%  "

$ws.Range("D36").Value = "'0.693"
$ws.Range("E36").Value = "'  +0.43%  "

$ws.Range("D37").Value = "'91.31"
$ws.Range("E37").Value = "'  -4.64%  "

$ws.Range("E38").Value = "'  +4.65%  "

$ws.Range("D39").Value = "'1.320.68"
$ws.Range("E39").Value = "'  -2.10%  "

$ws.Range("E40").Value = "'  -2.19%  "

$ws.Range("E41").Value = "'  -0.22%  "

$ws.Range("D42").Value = "'0.957"
$ws.Range("E42").Value = "'  -5.63%  "

$ws.Range("D43").Value = "'14.25"
$ws.Range("E43").Value = "'  -8.62%  "

$ws.Range("E44").Value = "'  -3.47%  "

$ws.Range("E45").Value = "'  -10.27%  "

$ws.Range("D46").Value = "'6.21"
$ws.Range("E46").Value = "'  -1.29%  "

$ws.Range("E47").Value = "'  -1.29%  "

$ws.Range("D48").Value = "'1.994.01"
$ws.Range("E48").Value = "'  -0.71%  "

$ws.Range("D50").Value = "'0.0663"
$ws.Range("E50").Value = "'  +4.89%  "

$ws.Range("D51").Value = "'97.90"
